$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1908.1428
$ws.Range("I40").Value = 760
$ws.Range("J40").Value = 2099.5
$ws.Range("K40").Value = 760
$ws.Range("L40").Value = 2099.5
$ws.Range("M40").Value = -585
$ws.Range("N40").Value = -2449.5

$ws.Range("H64").Value = 3520
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3520
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3520
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4016

$ws.Range("H67").Value = 3520
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3520
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3520
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5236

$ws.Range("H76").Value = 2781047.8
$ws.Range("I76").Value = 3457.1428
$ws.Range("K76").Value = 3457.1428
$ws.Range("M76").Value = -3142.1428

$ws.Range("H79").Value = 2781047.8
$ws.Range("I79").Value = 3457.1428
$ws.Range("K79").Value = 3457.1428
$ws.Range("M79").Value = -2365.1428

$ws.Range("H129").Value = 197076.47
$ws.Range("J129").Value = 213819.36
$ws.Range("L129").Value = 641458.08
$ws.Range("N129").Value = -651458.08

$ws.Range("H137").Value = 1181.9375
$ws.Range("I137").Value = 1126.9546
$ws.Range("K137").Value = 3380.8638
$ws.Range("M137").Value = -830.8638000000001

$ws.Range("H138").Value = 2098.261
$ws.Range("I138").Value = 1271.7727
$ws.Range("J138").Value = 2485.1277
$ws.Range("K138").Value = 3815.3181
$ws.Range("L138").Value = 7455.3831
$ws.Range("M138").Value = 1324.6819
$ws.Range("N138").Value = -17735.3831

$ws.Range("H141").Value = 3031.6
$ws.Range("I141").Value = 2830.4443
$ws.Range("J141").Value = 3333.3333
$ws.Range("K141").Value = 8491.332900000001
$ws.Range("L141").Value = 9999.999899999999
$ws.Range("M141").Value = -3311.332900000001
$ws.Range("N141").Value = -20359.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6204.6094
$ws.Range("I32").Value = 4955.609
$ws.Range("J32").Value = 10992.444
$ws.Range("K32").Value = 4955.609
$ws.Range("L32").Value = 10992.444
$ws.Range("M32").Value = -4668.609
$ws.Range("N32").Value = -11566.444

$ws.Range("H63").Value = 3126815
$ws.Range("I63").Value = 2016.7778
$ws.Range("K63").Value = 2016.7778
$ws.Range("M63").Value = -1330.7778

$ws.Range("H66").Value = 3126815
$ws.Range("I66").Value = 2016.7778
$ws.Range("K66").Value = 10083.889
$ws.Range("M66").Value = -6651.889000000001

$ws.Range("H74").Value = 31251796
$ws.Range("I74").Value = 62500700
$ws.Range("J74").Value = 2894.5
$ws.Range("K74").Value = 62500700
$ws.Range("L74").Value = 2894.5
$ws.Range("M74").Value = -62499826
$ws.Range("N74").Value = -4642.5

$ws.Range("H77").Value = 31251796
$ws.Range("I77").Value = 62500700
$ws.Range("J77").Value = 2894.5
$ws.Range("K77").Value = 312503500
$ws.Range("L77").Value = 14472.5
$ws.Range("M77").Value = -312499132
$ws.Range("N77").Value = -23208.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 14907
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 14907
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H105").Value = 4168526
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 5001831
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 5001831
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -5005325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3616.641
$ws.Range("I31").Value = 1892.2727
$ws.Range("K31").Value = 1892.2727
$ws.Range("M31").Value = -1597.2727

$ws.Range("H34").Value = 3616.641
$ws.Range("I34").Value = 1892.2727
$ws.Range("K34").Value = 1892.2727
$ws.Range("M34").Value = -1690.2727

$ws.Range("H62").Value = 5983.3335
$ws.Range("I62").Value = 5975
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 5975
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -5351
$ws.Range("N62").Value = -7248

$ws.Range("H65").Value = 5983.3335
$ws.Range("I65").Value = 5975
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 29875
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -26755
$ws.Range("N65").Value = -36240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 686.24
$ws.Range("J131").Value = 713.45056
$ws.Range("L131").Value = 2140.35168
$ws.Range("N131").Value = -12220.35168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4475165
$ws.Range("I70").Value = 4718
$ws.Range("J70").Value = 7828000
$ws.Range("K70").Value = 4718
$ws.Range("L70").Value = 7828000
$ws.Range("M70").Value = -4448
$ws.Range("N70").Value = -7828540

$ws.Range("H73").Value = 4475165
$ws.Range("I73").Value = 4718
$ws.Range("J73").Value = 7828000
$ws.Range("K73").Value = 4718
$ws.Range("L73").Value = 7828000
$ws.Range("M73").Value = -3782
$ws.Range("N73").Value = -7829872

$ws.Range("H80").Value = 3514.7407
$ws.Range("I80").Value = 3091.0833
$ws.Range("J80").Value = 3853.6667
$ws.Range("K80").Value = 3091.0833
$ws.Range("L80").Value = 3853.6667
$ws.Range("M80").Value = -2093.0833
$ws.Range("N80").Value = -5849.6667

$ws.Range("H83").Value = 3514.7407
$ws.Range("I83").Value = 3091.0833
$ws.Range("J83").Value = 3853.6667
$ws.Range("K83").Value = 15455.4165
$ws.Range("L83").Value = 19268.3335
$ws.Range("M83").Value = -10463.4165
$ws.Range("N83").Value = -29252.3335

$ws.Range("H102").Value = 1575.5483
$ws.Range("I102").Value = 1284.9286
$ws.Range("K102").Value = 1284.9286
$ws.Range("M102").Value = 337.0714

$ws.Range("H132").Value = 21824.143
$ws.Range("I132").Value = 4075.4285
$ws.Range("J132").Value = 75070.28999999999
$ws.Range("K132").Value = 12226.2855
$ws.Range("L132").Value = 225210.87
$ws.Range("M132").Value = -9696.2855
$ws.Range("N132").Value = -230270.87

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 403308.03
$ws.Range("I132").Value = 483050.53
$ws.Range("J132").Value = 4595.6
$ws.Range("K132").Value = 1449151.59
$ws.Range("L132").Value = 13786.8
$ws.Range("M132").Value = -1446621.59
$ws.Range("N132").Value = -18846.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 5003000
$ws.Range("J18").Value = 5003000
$ws.Range("L18").Value = 5003000
$ws.Range("N18").Value = -5003346

$ws.Range("J58").Value = 18200
$ws.Range("L58").Value = 18200
$ws.Range("N58").Value = -18816

$ws.Range("H86").Value = 25750
$ws.Range("J86").Value = 25750
$ws.Range("L86").Value = 25750
$ws.Range("N86").Value = -27996

$ws.Range("H89").Value = 25750
$ws.Range("J89").Value = 25750
$ws.Range("L89").Value = 128750
$ws.Range("N89").Value = -139982

$ws.Range("H96").Value = 1886.6666
$ws.Range("I96").Value = 1775
$ws.Range("J96").Value = 1976
$ws.Range("K96").Value = 1775
$ws.Range("L96").Value = 1976
$ws.Range("M96").Value = -402
$ws.Range("N96").Value = -4722

$ws.Range("H100").Value = 295.33334
$ws.Range("I100").Value = 267.5
$ws.Range("K100").Value = 535
$ws.Range("M100").Value = 6

$ws.Range("H113").Value = 1658.8667
$ws.Range("I113").Value = 1658.8667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4976.6001
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2806.6001
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 1344.2195
$ws.Range("I132").Value = 1075.5358
$ws.Range("J132").Value = 1922.9231
$ws.Range("K132").Value = 3226.6074
$ws.Range("L132").Value = 5768.7693
$ws.Range("M132").Value = -696.6074000000003
